# Auto-generated Excel COM-interop edit script
# Updates the cryptos price/volume table to the latest scrape values,
# including a couple of row swaps (rank re-ordering) for TRON/Litecoin
# and RenderToken/TheSandbox.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates (coin names, links, percentage strings) ---
# These are not numeric-looking, so a normal .Value2 assignment keeps
# them stored as text, matching the inline-string cells in the sheet.
$ws.Range("D2").Value2 = '27.367.85'
$ws.Range("E2").Value2 = '  +2.64%  '
$ws.Range("D3").Value2 = '1.822.20'
$ws.Range("E3").Value2 = '  +1.66%  '
$ws.Range("E4").Value2 = '  -0.38%  '
$ws.Range("E5").Value2 = '  +1.76%  '
$ws.Range("E6").Value2 = '  -0.22%  '
$ws.Range("E7").Value2 = '  +5.67%  '
$ws.Range("E8").Value2 = '  +2.98%  '
$ws.Range("E9").Value2 = '  +1.61%  '
$ws.Range("E10").Value2 = '  +2.09%  '
$ws.Range("E11").Value2 = '  +0.71%  '
$ws.Range("D12").Value2 = '1.825.53'
$ws.Range("E12").Value2 = '  +1.80%  '
$ws.Range("E13").Value2 = '  +1.35%  '
$ws.Range("E14").Value2 = '  +3.27%  '
$ws.Range("B15").Value2 = 'TRON'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("E15").Value2 = '  +0.70%  '
$ws.Range("B16").Value2 = 'Litecoin'
$ws.Range("C16").Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("E16").Value2 = '  +1.38%  '
$ws.Range("E17").Value2 = '  -0.36%  '
$ws.Range("E18").Value2 = '  +1.40%  '
$ws.Range("E19").Value2 = '  -0.12%  '
$ws.Range("E20").Value2 = '  +1.52%  '
$ws.Range("D21").Value2 = '27.371.22'
$ws.Range("E21").Value2 = '  +2.56%  '
$ws.Range("E22").Value2 = '  +3.25%  '
$ws.Range("E23").Value2 = '  +1.31%  '
$ws.Range("D24").Value2 = '2.052.48'
$ws.Range("E24").Value2 = '  +1.78%  '
$ws.Range("E25").Value2 = '  -1.66%  '
$ws.Range("E26").Value2 = '  -0.03%  '
$ws.Range("E27").Value2 = '  +4.16%  '
$ws.Range("E28").Value2 = '  +1.54%  '
$ws.Range("E29").Value2 = '  +3.11%  '
$ws.Range("E30").Value2 = '  +0.25%  '
$ws.Range("E31").Value2 = '  +1.68%  '
$ws.Range("E32").Value2 = '  +6.98%  '
$ws.Range("E33").Value2 = '  +3.89%  '
$ws.Range("E34").Value2 = '  +2.74%  '
$ws.Range("E35").Value2 = '  +0.71%  '
$ws.Range("E36").Value2 = '  -0.31%  '
$ws.Range("E37").Value2 = '  +1.22%  '
$ws.Range("E38").Value2 = '  +1.61%  '
$ws.Range("E39").Value2 = '  +2.14%  '
$ws.Range("E40").Value2 = '  +5.17%  '
$ws.Range("B41").Value2 = 'RenderToken'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E41").Value2 = '  +21.99%  '
$ws.Range("B42").Value2 = 'TheSandbox'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("E42").Value2 = '  +1.93%  '
$ws.Range("E43").Value2 = '  +3.41%  '
$ws.Range("E44").Value2 = '  +1.15%  '
$ws.Range("E45").Value2 = '  +2.95%  '
$ws.Range("E46").Value2 = '  +0.99%  '
$ws.Range("E47").Value2 = '  +1.52%  '
$ws.Range("E48").Value2 = '  +0.95%  '
$ws.Range("E49").Value2 = '  +1.43%  '
$ws.Range("E50").Value2 = '  -0.18%  '
$ws.Range("E51").Value2 = '  +0.88%  '

# --- Price column (D) updates that LOOK like numbers -------------
# The source sheet stores these as literal text (e.g. "1.0000",
# "0.9989", "0.000008776") rather than numeric values, so a plain
# assignment would let Excel's smart-parsing coerce them into
# doubles and mangle formatting (trailing zeros, sci-notation, etc).
# Temporarily force the Text number format, write the literal string,
# then restore the cell style so no stray formatting is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '0.9989'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '313.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4659'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.3773'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.07441'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.8722'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '20.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '6.685'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '5.414'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '0.07103'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '92.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '0.9998'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.000008776'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '1.0000'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '14.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '10.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '1.940'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '151.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '2.261'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '18.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '5.302'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '117.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '0.08899'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '0.7837'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '1.185'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.536'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '2.925'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.9991'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '1.098'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.01974'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.05265'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '7.295'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '2.389'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.5309'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '2.904'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.1691'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '8.635'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.5058'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '10.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '105.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '1.678'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.9991'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '0.06339'
$ws.Range("D51").Style = "Normal"
